# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ /
# LeveProfitHQ figures (columns H-N) across the job-sheet tables, refreshed
# from the latest market-board pull (scheduled runner sync).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1451.4706
$ws.Range("J17").Value = 1451.4706
$ws.Range("L17").Value = 4354.4118
$ws.Range("N17").Value = -4690.4118
$ws.Range("H18").Value = 12196.4
$ws.Range("I18").Value = 7000
$ws.Range("J18").Value = 12773.777
$ws.Range("K18").Value = 7000
$ws.Range("L18").Value = 12773.777
$ws.Range("M18").Value = -6716
$ws.Range("N18").Value = -13341.777
$ws.Range("H28").Value = 327.63635
$ws.Range("I28").Value = 259.42856
$ws.Range("J28").Value = 447
$ws.Range("K28").Value = 259.42856
$ws.Range("L28").Value = 447
$ws.Range("M28").Value = 225.57144
$ws.Range("N28").Value = -1417
$ws.Range("H33").Value = 101.70588
$ws.Range("J33").Value = 113.63636
$ws.Range("L33").Value = 113.63636
$ws.Range("N33").Value = -571.63636
$ws.Range("H64").Value = 2999.6667
$ws.Range("I64").Value = 2910.4
$ws.Range("J64").Value = 3111.25
$ws.Range("K64").Value = 2910.4
$ws.Range("L64").Value = 3111.25
$ws.Range("M64").Value = -2662.4
$ws.Range("N64").Value = -3607.25
$ws.Range("H67").Value = 2999.6667
$ws.Range("I67").Value = 2910.4
$ws.Range("J67").Value = 3111.25
$ws.Range("K67").Value = 2910.4
$ws.Range("L67").Value = 3111.25
$ws.Range("M67").Value = -2052.4
$ws.Range("N67").Value = -4827.25
$ws.Range("H74").Value = 2999.5
$ws.Range("I74").Value = 2999.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2999.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -2063.5
$ws.Range("N74").Value = $null
$ws.Range("H77").Value = 2999.5
$ws.Range("I77").Value = 2999.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 14997.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -10317.5
$ws.Range("N77").Value = $null
$ws.Range("H113").Value = 18850
$ws.Range("I113").Value = 20290.908
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 20290.908
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -17036.908
$ws.Range("N113").Value = -9508

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3311.9849
$ws.Range("I32").Value = 2185.0566
$ws.Range("K32").Value = 2185.0566
$ws.Range("M32").Value = -1898.0566
$ws.Range("H61").Value = 4538.6665
$ws.Range("I61").Value = 2362.4
$ws.Range("J61").Value = 7259
$ws.Range("K61").Value = 2362.4
$ws.Range("L61").Value = 7259
$ws.Range("M61").Value = -2150.4
$ws.Range("N61").Value = -7683
$ws.Range("H74").Value = 1240.0869
$ws.Range("I74").Value = 1127
$ws.Range("K74").Value = 1127
$ws.Range("M74").Value = -253
$ws.Range("H77").Value = 1240.0869
$ws.Range("I77").Value = 1127
$ws.Range("K77").Value = 5635
$ws.Range("M77").Value = -1267
$ws.Range("H132").Value = 1658.125
$ws.Range("I132").Value = 1285.9412
$ws.Range("K132").Value = 3857.8236
$ws.Range("M132").Value = -1327.8236
$ws.Range("H136").Value = 4538.6665
$ws.Range("I136").Value = 2362.4
$ws.Range("J136").Value = 7259
$ws.Range("K136").Value = 7087.200000000001
$ws.Range("L136").Value = 21777
$ws.Range("M136").Value = -4537.200000000001
$ws.Range("N136").Value = -26877

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 950.8333
$ws.Range("I22").Value = 375
$ws.Range("K22").Value = 375
$ws.Range("M22").Value = -25
$ws.Range("H58").Value = 1360362.8
$ws.Range("I58").Value = 2175003.8
$ws.Range("J58").Value = 2627.75
$ws.Range("K58").Value = 2175003.8
$ws.Range("L58").Value = 2627.75
$ws.Range("M58").Value = -2174800.8
$ws.Range("N58").Value = -3033.75
$ws.Range("H60").Value = 11360.4
$ws.Range("J60").Value = 11360.4
$ws.Range("L60").Value = 11360.4
$ws.Range("N60").Value = -12382.4
$ws.Range("H105").Value = 2203.1667
$ws.Range("J105").Value = 2500
$ws.Range("L105").Value = 2500
$ws.Range("N105").Value = -5994
$ws.Range("H134").Value = 1882.1333
$ws.Range("I134").Value = 1683.5
$ws.Range("J134").Value = 3173.25
$ws.Range("K134").Value = 5050.5
$ws.Range("L134").Value = 9519.75
$ws.Range("M134").Value = -2515.5
$ws.Range("N134").Value = -14589.75
$ws.Range("H136").Value = 1360362.8
$ws.Range("I136").Value = 2175003.8
$ws.Range("J136").Value = 2627.75
$ws.Range("K136").Value = 6525011.399999999
$ws.Range("L136").Value = 7883.25
$ws.Range("M136").Value = -6522461.399999999
$ws.Range("N136").Value = -12983.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 5250
$ws.Range("I70").Value = 500
$ws.Range("J70").Value = 10000
$ws.Range("K70").Value = 1500
$ws.Range("L70").Value = 30000
$ws.Range("M70").Value = -1185
$ws.Range("N70").Value = -30630
$ws.Range("H73").Value = 5250
$ws.Range("I73").Value = 500
$ws.Range("J73").Value = 10000
$ws.Range("K73").Value = 1500
$ws.Range("L73").Value = 30000
$ws.Range("M73").Value = -408
$ws.Range("N73").Value = -32184
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = $null
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = $null
$ws.Range("H131").Value = 758.78
$ws.Range("J131").Value = 773.94684
$ws.Range("L131").Value = 2321.84052
$ws.Range("N131").Value = -12401.84052

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1874.75
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 2166.3333
$ws.Range("K80").Value = 1000
$ws.Range("L80").Value = 2166.3333
$ws.Range("M80").Value = -2
$ws.Range("N80").Value = -4162.3333
$ws.Range("H83").Value = 1874.75
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 2166.3333
$ws.Range("K83").Value = 5000
$ws.Range("L83").Value = 10831.6665
$ws.Range("M83").Value = -8
$ws.Range("N83").Value = -20815.6665
$ws.Range("H122").Value = 1980.5883
$ws.Range("J122").Value = 2642.5715
$ws.Range("L122").Value = 7927.7145
$ws.Range("N122").Value = -12827.7145
$ws.Range("H126").Value = 2461697.8
$ws.Range("I126").Value = 2528138.5
$ws.Range("K126").Value = 7584415.5
$ws.Range("M126").Value = -7581945.5
$ws.Range("H132").Value = 2140792.8
$ws.Range("J132").Value = 3876.5
$ws.Range("L132").Value = 11629.5
$ws.Range("N132").Value = -16689.5
$ws.Range("H136").Value = 8360.857
$ws.Range("J136").Value = 8360.857
$ws.Range("L136").Value = 25082.571
$ws.Range("N136").Value = -30182.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3210
$ws.Range("I7").Value = 2368.1
$ws.Range("J7").Value = 4893.8
$ws.Range("K7").Value = 2368.1
$ws.Range("L7").Value = 4893.8
$ws.Range("M7").Value = -2256.1
$ws.Range("N7").Value = -5117.8
$ws.Range("H126").Value = 3210
$ws.Range("I126").Value = 2368.1
$ws.Range("J126").Value = 4893.8
$ws.Range("K126").Value = 7104.299999999999
$ws.Range("L126").Value = 14681.4
$ws.Range("M126").Value = -4634.299999999999
$ws.Range("N126").Value = -19621.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2155.2
$ws.Range("I81").Value = 2194.25
$ws.Range("J81").Value = 1999
$ws.Range("K81").Value = 4388.5
$ws.Range("L81").Value = 3998
$ws.Range("M81").Value = -3327.5
$ws.Range("N81").Value = -6120
$ws.Range("H84").Value = 2155.2
$ws.Range("I84").Value = 2194.25
$ws.Range("J84").Value = 1999
$ws.Range("K84").Value = 21942.5
$ws.Range("L84").Value = 19990
$ws.Range("M84").Value = -16638.5
$ws.Range("N84").Value = -30598

